$d = $word.ActiveDocument

# 1. Color the title "EL PLASTICO Y LA SALUD" red (FF0000)
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Font.Color = 255

# 2. Merge the split runs "A" / " " / "pesar de ser..." into a single run
$d.Content.Find.Execute("A pesar de ser", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "A pesar de ser", 2)

# 3. Remove the stale <w:lastRenderedPageBreak/> marker preceding the drawing.
#    The Word OM has no direct property for this rendering-cache artifact, so we
#    insert a clean copy of the drawing run (identical XML, minus the marker)
#    immediately before the existing one, then delete the original drawing shape
#    (which removes its whole run, marker included) leaving only the clean run.
$lastPara = $d.Paragraphs.Last
$insPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$drawingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251658240" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="03BB6B0C" wp14:editId="5BAA3469"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="margin"><wp:posOffset>-1080135</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>0</wp:posOffset></wp:positionV><wp:extent cx="7543800" cy="7305675"/><wp:effectExtent l="0" t="0" r="0" b="9525"/><wp:wrapSquare wrapText="bothSides"/><wp:docPr id="1" name="Imagen 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1" name="Imagen 1"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId4"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="7543800" cy="7305675"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></w:r></w:p>'
$insPoint.InsertXML($drawingXml)
$oldShape = $d.Shapes.Item(2)
$oldShape.Delete()
